{"js": "// The character \"Danniel\" was renamed to \"Dimitri\" everywhere he is\n// mentioned in the storyboard text: the title card (\"corporal\n// 'Danniel'\") and the two later narrative mentions (\"maar Danniel\n// wist...\" / \"Danniel probeert nu...\"). Find every occurrence of the\n// old name in the document body and replace it with the new one.\n\nconst body = context.document.body;\n\nconst nameHits = body.search(\"Danniel\", { matchCase: true, matchWholeWord: true });\nnameHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of nameHits.items) {\n  hit.insertText(\"Dimitri\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// The sentence \"Hij kwam van ruimtestation 'Omega'. Hij leefde een\n// normaal leven ... totdat het werd aangevallen\" was also touched while\n// the author was editing that paragraph (no wording changed, Word just\n// re-flowed a few adjacent sentences into a single run). Re-assert the\n// same text in place so the run boundaries match.\nconst reflowed = \". Hij leefde een normaal leven en deed de taken binnen de ruimtestation totdat het\";\nconst reflowHits = body.search(reflowed, { matchCase: true });\nreflowHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of reflowHits.items) {\n  hit.insertText(reflowed, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The character \"Danniel\" was renamed to \"Dimitri\" everywhere he is\n# mentioned in the storyboard text (title card, and the two later\n# narrative mentions). Find every occurrence of the old name in the\n# document and replace it with the new one.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Danniel\"\n$find.Replacement.Text = \"Dimitri\"\n\n# wdFindContinue = 1 (search whole story, no \"replace past end\" prompt),\n# wdReplaceAll = 2 -> replace every match in one pass.\n$find.Execute(\n    \"Danniel\",  # FindText\n    $true,      # MatchCase\n    $true,      # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    \"Dimitri\",  # ReplaceWith\n    2           # Replace (wdReplaceAll)\n)\n\nWrite-Output \"done\"\n"}
